$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '57.186.54'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +2.42%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.258.21'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.32%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '398.07'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.25%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '108.98'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.61%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.582'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +4.83%  '
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("E9").Value = '  -0.18%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.40'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.28%  '
$ws.Range("E11").Value = '  +5.20%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.774.02'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.45%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '8.29'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.55%  '
$ws.Range("E15").Value = '  -0.46%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.260.63'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.10%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.04'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.61%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '11.10'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +4.05%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '57.028.98'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.40%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '3.32'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.48%  '
$ws.Range("E21").Value = '  +6.31%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '12.95'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.14%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '294.42'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -3.33%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '74.23'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.44%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.19'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.09%  '
$ws.Range("B26").Value = 'EthereumClassic'
$ws.Range("C26").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '28.16'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.42%  '
$ws.Range("B27").Value = 'Filecoin'
$ws.Range("C27").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.92'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -4.00%  '
$ws.Range("E28").Value = '  +0.56%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.43'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.35%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.170'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.55%  '
$ws.Range("E31").Value = '  +0.00%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.111'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.59%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '11.21'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.85%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '40.21'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +11.37%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0489'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.15%  '
$ws.Range("E36").Value = '  +1.18%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '51.31'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.26%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.00'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.06%  '
$ws.Range("E39").Value = '  -0.74%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.01'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.85%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '136.49'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.24%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.95'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.17%  '
$ws.Range("E44").Value = '  -2.62%  '
$ws.Range("E45").Value = '  -1.23%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '16.81'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.04%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '22.50'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.23%  '
$ws.Range("E48").Value = '  +5.57%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.147.99'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.44%  '
$ws.Range("E50").Value = '  -7.60%  '
$ws.Range("E51").Value = '  -6.62%  '
